$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.304.20"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.921.20"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - XRP
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.8136"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.89%  "

# Row 6 - BNB
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "244.21"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "

# Row 7 - USDC
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Cardano
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3268"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.37%  "

# Row 9 - Solana
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "27.27"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.61%  "

# Row 10 - Dogecoin
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07273"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.39%  "

# Row 11 - Polygon
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.7952"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.05%  "

# Row 12 - TRON
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08111"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.37%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.934.07"
$ws.Range("E13").Value = "  +1.28%  "

# Row 14 - Polkadot
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.411"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.35%  "

# Row 15 - Litecoin
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "94.16"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "30.330.70"

# Row 17 - Avalanche
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.28"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.30%  "

# Row 18 - Uniswap
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.080"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.69%  "

# Row 19 - BitcoinCash
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "251.28"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.31%  "

# Row 20 - ShibaInu
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007862"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.179.58"
$ws.Range("E21").Value = "  +0.90%  "

# Row 22 - Dai
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23 - Chainlink
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.039"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +17.63%  "

# Row 24 - BinanceUSD
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25 - Stellar
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1680"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +21.11%  "

# Row 26 - Cosmos
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.510"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.07%  "

# Row 27 - Monero
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "167.62"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "

# Row 28 - EthereumClassic
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.09"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "

# Row 29 - LidoDAOToken
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.162"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.45%  "

# Row 30 - Toncoin
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.372"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "

# Row 31 - PancakeSwap
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.550"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.37%  "

# Row 32 - Filecoin
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.354"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "

# Row 33 - Hedera
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05685"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.01%  "

# Row 34 - InternetComputer(DFINITY)
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.144"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.51%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +3.52%  "

# Row 36 - ImmutableX
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7468"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.97%  "

# Row 37 - Frax
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38 - HuobiToken
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.730"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

# Row 39 - VeChain
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01962"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "

# Row 40 - MXToken
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.818"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.27%  "

# Row 41 - TheSandbox
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.4509"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.23%  "

# Row 42 - Aave
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "74.77"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.30%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -2.33%  "

# Row 44 - TrustWalletToken
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.8560"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.20%  "

# Row 45 - RenderToken
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.931"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.07%  "

# Row 46 - Maker (-> PaxDollar)
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47 - PaxDollar (-> Maker)
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.038.96"
$ws.Range("E47").Value = "  +5.18%  "

# Row 48 - Quant
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "103.16"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.77%  "

# Row 49 - Aptos (-> SynthetixNetwork)
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.122"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +11.89%  "

# Row 50 - EnergySwap
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "9.993"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "

# Row 51 - SynthetixNetwork (-> Aptos)
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.653"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.42%  "
